$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143. This shifts existing rows 143..279 down to 144..280.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with a new data record (same as the former
# row 143 except for a new "Fecha" date).
$newRow = 143

$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value = "Maule"

$dCell = $ws.Cells.Item($newRow, 4)
$dCell.Value = Get-Date -Year 2022 -Month 4 -Day 18 -Hour 0 -Minute 0 -Second 0
$dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = 100112003
$ws.Cells.Item($newRow, 7).Value = "Ajo"
$ws.Cells.Item($newRow, 8).Value = "Chino"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 300
$ws.Cells.Item($newRow, 11).Value = 21000
$ws.Cells.Item($newRow, 12).Value = 21000
$ws.Cells.Item($newRow, 13).Value = 21000
$ws.Cells.Item($newRow, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item($newRow, 15).Value = "China"
$ws.Cells.Item($newRow, 16).Value = 2100
$ws.Cells.Item($newRow, 17).Value = 10
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
